$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '27.438.77'
$ws.Cells.Item(2, 5).NumberFormat = "@"
$ws.Cells.Item(2, 5).Value = '  +0.73%  '
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '1.719.71'
$ws.Cells.Item(3, 5).NumberFormat = "@"
$ws.Cells.Item(3, 5).Value = '  +0.48%  '
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '1.006'
$ws.Cells.Item(4, 5).NumberFormat = "@"
$ws.Cells.Item(4, 5).Value = '  +1.13%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '224.79'
$ws.Cells.Item(5, 5).NumberFormat = "@"
$ws.Cells.Item(5, 5).Value = '  -2.30%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '0.5332'
$ws.Cells.Item(6, 5).NumberFormat = "@"
$ws.Cells.Item(6, 5).Value = '  -1.94%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '1.006'
$ws.Cells.Item(7, 5).NumberFormat = "@"
$ws.Cells.Item(7, 5).Value = '  +1.06%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.2653'
$ws.Cells.Item(8, 5).NumberFormat = "@"
$ws.Cells.Item(8, 5).Value = '  -2.89%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.06582'
$ws.Cells.Item(9, 5).NumberFormat = "@"
$ws.Cells.Item(9, 5).Value = '  +1.12%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '21.23'
$ws.Cells.Item(10, 5).NumberFormat = "@"
$ws.Cells.Item(10, 5).Value = '  -1.89%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.07668'
$ws.Cells.Item(11, 5).NumberFormat = "@"
$ws.Cells.Item(11, 5).Value = '  -1.19%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '4.600'
$ws.Cells.Item(12, 5).NumberFormat = "@"
$ws.Cells.Item(12, 5).Value = '  -2.92%  '
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '1.711.57'
$ws.Cells.Item(13, 5).NumberFormat = "@"
$ws.Cells.Item(13, 5).Value = '  -0.37%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '1.959.39'
$ws.Cells.Item(14, 5).NumberFormat = "@"
$ws.Cells.Item(14, 5).Value = '  +1.14%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '0.5772'
$ws.Cells.Item(15, 5).NumberFormat = "@"
$ws.Cells.Item(15, 5).Value = '  -3.20%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '0.0₅8244'
$ws.Cells.Item(16, 5).NumberFormat = "@"
$ws.Cells.Item(16, 5).Value = '  -1.41%  '
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '67.71'
$ws.Cells.Item(17, 5).NumberFormat = "@"
$ws.Cells.Item(17, 5).Value = '  -0.68%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '27.454.43'
$ws.Cells.Item(18, 5).NumberFormat = "@"
$ws.Cells.Item(18, 5).Value = '  +1.37%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '217.70'
$ws.Cells.Item(19, 5).NumberFormat = "@"
$ws.Cells.Item(19, 5).Value = '  +3.45%  '
$ws.Cells.Item(20, 5).NumberFormat = "@"
$ws.Cells.Item(20, 5).Value = '  +0.96%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '4.717'
$ws.Cells.Item(21, 5).NumberFormat = "@"
$ws.Cells.Item(21, 5).Value = '  -1.75%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '10.53'
$ws.Cells.Item(22, 5).NumberFormat = "@"
$ws.Cells.Item(22, 5).Value = '  -4.31%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '5.976'
$ws.Cells.Item(23, 5).NumberFormat = "@"
$ws.Cells.Item(23, 5).Value = '  -3.75%  '
$ws.Cells.Item(24, 5).NumberFormat = "@"
$ws.Cells.Item(24, 5).Value = '  +1.08%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '143.06'
$ws.Cells.Item(25, 5).NumberFormat = "@"
$ws.Cells.Item(25, 5).Value = '  -2.59%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '1.739'
$ws.Cells.Item(26, 5).NumberFormat = "@"
$ws.Cells.Item(26, 5).Value = '  +10.42%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '0.1226'
$ws.Cells.Item(27, 5).NumberFormat = "@"
$ws.Cells.Item(27, 5).Value = '  -2.80%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '7.301'
$ws.Cells.Item(28, 5).NumberFormat = "@"
$ws.Cells.Item(28, 5).Value = '  -1.55%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '16.42'
$ws.Cells.Item(29, 5).NumberFormat = "@"
$ws.Cells.Item(29, 5).Value = '  -1.95%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '0.05430'
$ws.Cells.Item(30, 5).NumberFormat = "@"
$ws.Cells.Item(30, 5).Value = '  -4.72%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '1.297'
$ws.Cells.Item(31, 5).NumberFormat = "@"
$ws.Cells.Item(31, 5).Value = '  -1.84%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '3.523'
$ws.Cells.Item(32, 5).NumberFormat = "@"
$ws.Cells.Item(32, 5).Value = '  -3.29%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '3.426'
$ws.Cells.Item(33, 5).NumberFormat = "@"
$ws.Cells.Item(33, 5).Value = '  -2.68%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '1.647'
$ws.Cells.Item(34, 5).NumberFormat = "@"
$ws.Cells.Item(34, 5).Value = '  +0.47%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '2.879'
$ws.Cells.Item(35, 5).NumberFormat = "@"
$ws.Cells.Item(35, 5).Value = '  +0.69%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.9534'
$ws.Cells.Item(36, 5).NumberFormat = "@"
$ws.Cells.Item(36, 5).Value = '  -2.94%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '2.427'
$ws.Cells.Item(37, 5).NumberFormat = "@"
$ws.Cells.Item(37, 5).Value = '  +0.81%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.5892'
$ws.Cells.Item(38, 5).NumberFormat = "@"
$ws.Cells.Item(38, 5).Value = '  +0.54%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.01639'
$ws.Cells.Item(39, 5).NumberFormat = "@"
$ws.Cells.Item(39, 5).Value = '  -0.24%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '5.894'
$ws.Cells.Item(40, 5).NumberFormat = "@"
$ws.Cells.Item(40, 5).Value = '  -0.98%  '
$ws.Cells.Item(41, 2).NumberFormat = "@"
$ws.Cells.Item(41, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(41, 3).NumberFormat = "@"
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.8445'
$ws.Cells.Item(41, 5).NumberFormat = "@"
$ws.Cells.Item(41, 5).Value = '  -0.04%  '
$ws.Cells.Item(42, 2).NumberFormat = "@"
$ws.Cells.Item(42, 2).Value = 'PaxDollar'
$ws.Cells.Item(42, 3).NumberFormat = "@"
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '1.005'
$ws.Cells.Item(42, 5).NumberFormat = "@"
$ws.Cells.Item(42, 5).Value = '  +1.15%  '
$ws.Cells.Item(43, 2).NumberFormat = "@"
$ws.Cells.Item(43, 2).Value = 'Maker'
$ws.Cells.Item(43, 3).NumberFormat = "@"
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '1.043.91'
$ws.Cells.Item(43, 5).NumberFormat = "@"
$ws.Cells.Item(43, 5).Value = '  -1.98%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '101.13'
$ws.Cells.Item(44, 5).NumberFormat = "@"
$ws.Cells.Item(44, 5).Value = '  -1.48%  '
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '1.866.29'
$ws.Cells.Item(45, 5).NumberFormat = "@"
$ws.Cells.Item(45, 5).Value = '  +1.11%  '
$ws.Cells.Item(46, 5).NumberFormat = "@"
$ws.Cells.Item(46, 5).Value = '  +10.23%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '58.41'
$ws.Cells.Item(47, 5).NumberFormat = "@"
$ws.Cells.Item(47, 5).Value = '  -2.38%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '0.4512'
$ws.Cells.Item(48, 5).NumberFormat = "@"
$ws.Cells.Item(48, 5).Value = '  +4.04%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '8.188'
$ws.Cells.Item(49, 5).NumberFormat = "@"
$ws.Cells.Item(49, 5).Value = '  +1.10%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '1.006'
$ws.Cells.Item(50, 5).NumberFormat = "@"
$ws.Cells.Item(50, 5).Value = '  +1.03%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.06587'
$ws.Cells.Item(51, 5).NumberFormat = "@"
$ws.Cells.Item(51, 5).Value = '  +14.61%  '
